$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price/Volume columns being updated keep their original text
# representation (e.g. "519.92", "0.100", "0.0520") instead of being
# auto-converted to numbers by Excel, by forcing Text format first.
$cells = @('D2', 'E2', 'D3', 'E3', 'E4', 'D5', 'E5', 'D6', 'E6', 'D8', 'E8', 'D9', 'E9', 'E10', 'D11', 'E11', 'D12', 'E12', 'E13', 'D14', 'E14', 'D15', 'E15', 'D16', 'E16', 'E17', 'D18', 'E18', 'D19', 'E19', 'D20', 'E20', 'D21', 'E21', 'D22', 'E22', 'D23', 'E23', 'D24', 'E24', 'E25', 'D26', 'E26', 'D27', 'E27', 'D28', 'E28', 'D29', 'E29', 'D30', 'E30', 'E31', 'D32', 'E32', 'E33', 'D34', 'E34', 'D35', 'E35', 'D36', 'E36', 'D37', 'E37', 'D38', 'E38', 'D39', 'E39', 'D40', 'E40', 'E41', 'D42', 'E42', 'E43', 'D44', 'E44', 'D45', 'E45', 'D46', 'E46', 'D47', 'E47', 'D48', 'E48', 'D49', 'E49', 'D50', 'E50', 'D51', 'E51')
foreach ($cellRef in $cells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

$ws.Range('D2').Value = '57.514.54'
$ws.Range('E2').Value = '  -2.13%  '
$ws.Range('D3').Value = '2.561.11'
$ws.Range('E3').Value = '  -3.63%  '
$ws.Range('E4').Value = '  +0.03%  '
$ws.Range('D5').Value = '519.92'
$ws.Range('E5').Value = '  -0.74%  '
$ws.Range('D6').Value = '143.09'
$ws.Range('E6').Value = '  -0.61%  '
$ws.Range('D8').Value = '0.561'
$ws.Range('E8').Value = '  -1.41%  '
$ws.Range('D9').Value = '2.574.92'
$ws.Range('E9').Value = '  -3.38%  '
$ws.Range('E10').Value = '  -4.01%  '
$ws.Range('D11').Value = '0.100'
$ws.Range('E11').Value = '  -2.36%  '
$ws.Range('D12').Value = '0.325'
$ws.Range('E12').Value = '  -2.92%  '
$ws.Range('E13').Value = '  -0.25%  '
$ws.Range('D14').Value = '3.014.42'
$ws.Range('E14').Value = '  -3.60%  '
$ws.Range('D15').Value = '57.477.12'
$ws.Range('E15').Value = '  -2.18%  '
$ws.Range('D16').Value = '20.11'
$ws.Range('E16').Value = '  -4.08%  '
$ws.Range('E17').Value = '  -2.51%  '
$ws.Range('D18').Value = '2.576.13'
$ws.Range('E18').Value = '  -3.22%  '
$ws.Range('D19').Value = '334.92'
$ws.Range('E19').Value = '  -1.18%  '
$ws.Range('D20').Value = '4.28'
$ws.Range('E20').Value = '  -2.49%  '
$ws.Range('D21').Value = '10.17'
$ws.Range('E21').Value = '  -2.19%  '
$ws.Range('D22').Value = '6.20'
$ws.Range('E22').Value = '  -3.02%  '
$ws.Range('D23').Value = '0.999'
$ws.Range('E23').Value = '  +0.05%  '
$ws.Range('D24').Value = '65.16'
$ws.Range('E24').Value = '  +1.71%  '
$ws.Range('E25').Value = '  -0.37%  '
$ws.Range('D26').Value = '0.402'
$ws.Range('E26').Value = '  -4.97%  '
$ws.Range('D27').Value = '0.998'
$ws.Range('E27').Value = '  -0.10%  '
$ws.Range('D28').Value = '2.680.07'
$ws.Range('E28').Value = '  -3.69%  '
$ws.Range('D29').Value = '6.93'
$ws.Range('E29').Value = '  -2.92%  '
$ws.Range('D30').Value = '0.0₃0747'
$ws.Range('E30').Value = '  -7.07%  '
$ws.Range('E31').Value = '  -0.01%  '
$ws.Range('D32').Value = '6.23'
$ws.Range('E32').Value = '  -6.89%  '
$ws.Range('E33').Value = '  -0.96%  '
$ws.Range('D34').Value = '18.58'
$ws.Range('E34').Value = '  -1.57%  '
$ws.Range('D35').Value = '148.82'
$ws.Range('E35').Value = '  -1.46%  '
$ws.Range('D36').Value = '4.02'
$ws.Range('E36').Value = '  -3.19%  '
$ws.Range('D37').Value = '1.13'
$ws.Range('E37').Value = '  -4.26%  '
$ws.Range('D38').Value = '0.842'
$ws.Range('E38').Value = '  -10.14%  '
$ws.Range('D39').Value = '36.11'
$ws.Range('E39').Value = '  -1.72%  '
$ws.Range('D40').Value = '0.828'
$ws.Range('E40').Value = '  -5.35%  '
$ws.Range('E41').Value = '  -1.24%  '
$ws.Range('D42').Value = '3.51'
$ws.Range('E42').Value = '  -2.25%  '
$ws.Range('E43').Value = '  -0.14%  '
$ws.Range('D44').Value = '268.97'
$ws.Range('E44').Value = '  -2.47%  '
$ws.Range('D45').Value = '0.0955'
$ws.Range('E45').Value = '  -1.23%  '
$ws.Range('D46').Value = '10.63'
$ws.Range('E46').Value = '  -0.24%  '
$ws.Range('D47').Value = '0.587'
$ws.Range('E47').Value = '  -3.86%  '
$ws.Range('D48').Value = '18.86'
$ws.Range('E48').Value = '  -4.29%  '
$ws.Range('D49').Value = '0.0520'
$ws.Range('E49').Value = '  -2.67%  '
$ws.Range('D50').Value = '1.965.36'
$ws.Range('E50').Value = '  -4.66%  '
$ws.Range('D51').Value = '4.54'
$ws.Range('E51').Value = '  -3.82%  '
